# Updates cryptos list values (price + volume change %) to match the
# latest scrape data, including the RenderToken / Binance-PegBSC-USD
# row swap (rows 30 and 31).
#
# For numeric-looking Price values that must remain plain text (matching
# the source inlineStr cells), a leading apostrophe forces Excel to keep
# them as text instead of auto-converting to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.831.03'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '3.369.01'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''572.37'
$ws.Range('D6').Value = '''136.67'
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '3.367.23'
$ws.Range('E8').Value = '  -0.66%  '
$ws.Range('E9').Value = '  -1.08%  '
$ws.Range('D10').Value = '''7.64'
$ws.Range('E10').Value = '  +2.15%  '
$ws.Range('E11').Value = '  -2.47%  '
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('D13').Value = '3.945.56'
$ws.Range('E13').Value = '  -0.62%  '
$ws.Range('E14').Value = '  +0.52%  '
$ws.Range('D15').Value = '''25.85'
$ws.Range('E15').Value = '  +1.90%  '
$ws.Range('D16').Value = '3.371.39'
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('E17').Value = '  -3.30%  '
$ws.Range('D18').Value = '60.989.79'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('E19').Value = '  -1.85%  '
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').Value = '''9.39'
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('D22').Value = '''373.90'
$ws.Range('E23').Value = '  -2.81%  '
$ws.Range('D24').Value = '3.511.68'
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  -1.74%  '
$ws.Range('D27').Value = '''71.02'
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').Value = '''1.78'
$ws.Range('E28').Value = '  +12.12%  '
$ws.Range('E29').Value = '  +9.36%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').Value = '''7.44'
$ws.Range('E31').Value = '  -2.78%  '
$ws.Range('E32').Value = '  -2.12%  '
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  +0.74%  '
$ws.Range('D36').Value = '''5.17'
$ws.Range('E36').Value = '  -4.23%  '
$ws.Range('D37').Value = '''6.85'
$ws.Range('E37').Value = '  -1.27%  '
$ws.Range('E38').Value = '  -1.18%  '
$ws.Range('D39').Value = '''164.26'
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('D40').Value = '''0.0761'
$ws.Range('E40').Value = '  -3.21%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').Value = '''0.774'
$ws.Range('E42').Value = '  -1.05%  '
$ws.Range('D43').Value = '''41.55'
$ws.Range('E43').Value = '  -0.47%  '
$ws.Range('E44').Value = '  -4.65%  '
$ws.Range('E45').Value = '  -1.07%  '
$ws.Range('E46').Value = '  -2.91%  '
$ws.Range('D47').Value = '''24.15'
$ws.Range('E47').Value = '  -1.93%  '
$ws.Range('D48').Value = '2.454.61'
$ws.Range('E48').Value = '  +4.11%  '
$ws.Range('E49').Value = '  -2.12%  '
$ws.Range('D50').Value = '''22.91'
$ws.Range('E50').Value = '  -1.73%  '
$ws.Range('D51').Value = '''2.42'
$ws.Range('E51').Value = '  +3.99%  '
